$wb = $excel.ActiveWorkbook

# --- Sheet: ccs_retrofits (sheet2 in package, 2nd sheet / "ccs_retrofits") ---
$ws2 = $wb.Worksheets.Item("ccs_retrofits")

$ws2.Range("D4").Value = 0.09302400000000004

$ws2.Range("D5").Value = 0.09302400000000004

$ws2.Range("D6").Value = 0.09302400000000004

$ws2.Range("D7").Value = 0.09302400000000004

$ws2.Range("D8").Value = 0.21236400000000008

$ws2.Range("D9").Value = 0.21236400000000008

$ws2.Range("D10").Value = 0.08208000000000001

$ws2.Range("D11").Value = 0.08208000000000001

$ws2.Range("D12").Value = 0.08208000000000001

$ws2.Range("D13").Value = 0.08208000000000001

$ws2.Range("D14").Value = 0.08208000000000001

$ws2.Range("D15").Value = 0.08208000000000001

$ws2.Range("D16").Value = 0.08208000000000001

$ws2.Range("D17").Value = 0.08208000000000001

$ws2.Range("D18").Value = 0.21236400000000008

$ws2.Range("D19").Value = 0.21236400000000008

$ws2.Range("D20").Value = 0.22995000000000002
$ws2.Range("E20").Value = 2200
$ws2.Range("F20").Value = 27.7
$ws2.Range("G20").Value = 3.1
$ws2.Range("H20").Value = 0.74765

$ws2.Range("D21").Value = 0.22995000000000002
$ws2.Range("E21").Value = 2200
$ws2.Range("F21").Value = 27.7
$ws2.Range("G21").Value = 3.1
$ws2.Range("H21").Value = 0.74765

$ws2.Range("D22").Value = 0.22995000000000002
$ws2.Range("E22").Value = 2200
$ws2.Range("F22").Value = 27.7
$ws2.Range("G22").Value = 3.1
$ws2.Range("H22").Value = 0.74765

$ws2.Range("D23").Value = 0.22995000000000002
$ws2.Range("E23").Value = 2200
$ws2.Range("F23").Value = 27.7
$ws2.Range("G23").Value = 3.1
$ws2.Range("H23").Value = 0.74765

$ws2.Range("D24").Value = 0.08755200000000002

$ws2.Range("D25").Value = 0.08208000000000001

$ws2.Range("D26").Value = 0.08208000000000001

$ws2.Range("D27").Value = 0.08208000000000001

$ws2.Range("D28").Value = 0.08755200000000002

$ws2.Range("D29").Value = 0.08755200000000002

$ws2.Range("D30").Value = 0.08755200000000002

$ws2.Range("D31").Value = 0.08755200000000002

$ws2.Range("D32").Value = 0.08755200000000002

$ws2.Range("D33").Value = 0.08208000000000001

$ws2.Range("D34").Value = 0.08208000000000001

$ws2.Range("D35").Value = 0.08755200000000002

$ws2.Range("D36").Value = 0.08755200000000002

$ws2.Range("D37").Value = 0.08755200000000002

$ws2.Range("D38").Value = 0.08755200000000002

$ws2.Range("D39").Value = 0.08755200000000002

$ws2.Range("D40").Value = 0.22995000000000002
$ws2.Range("E40").Value = 2200
$ws2.Range("F40").Value = 27.7
$ws2.Range("G40").Value = 3.1
$ws2.Range("H40").Value = 0.74765

$ws2.Range("D41").Value = 0.22995000000000002
$ws2.Range("E41").Value = 2200
$ws2.Range("F41").Value = 27.7
$ws2.Range("G41").Value = 3.1
$ws2.Range("H41").Value = 0.74765

$ws2.Range("D42").Value = 0.22995000000000002
$ws2.Range("E42").Value = 2200
$ws2.Range("F42").Value = 27.7
$ws2.Range("G42").Value = 3.1
$ws2.Range("H42").Value = 0.74765

$ws2.Range("D43").Value = 0.23652
$ws2.Range("E43").Value = 2200
$ws2.Range("F43").Value = 27.7
$ws2.Range("G43").Value = 3.1
$ws2.Range("H43").Value = 0.74765

$ws2.Range("D44").Value = 0.23652
$ws2.Range("E44").Value = 2200
$ws2.Range("F44").Value = 27.7
$ws2.Range("G44").Value = 3.1
$ws2.Range("H44").Value = 0.74765

$ws2.Range("D45").Value = 0.23652
$ws2.Range("E45").Value = 2200
$ws2.Range("F45").Value = 27.7
$ws2.Range("G45").Value = 3.1
$ws2.Range("H45").Value = 0.74765

$ws2.Range("D46").Value = 0.21236400000000008

$ws2.Range("D47").Value = 0.21236400000000008

$ws2.Range("D48").Value = 0.21236400000000008

$ws2.Range("D49").Value = 0.21236400000000008

$ws2.Range("D50").Value = 0.21236400000000008

$ws2.Range("D51").Value = 0.21236400000000008

$ws2.Range("D52").Value = 0.22995000000000002
$ws2.Range("E52").Value = 2200
$ws2.Range("F52").Value = 27.7
$ws2.Range("G52").Value = 3.1
$ws2.Range("H52").Value = 0.74765

$ws2.Range("D53").Value = 0.22995000000000002
$ws2.Range("E53").Value = 2200
$ws2.Range("F53").Value = 27.7
$ws2.Range("G53").Value = 3.1
$ws2.Range("H53").Value = 0.74765

$ws2.Range("D54").Value = 0.22995000000000002
$ws2.Range("E54").Value = 2200
$ws2.Range("F54").Value = 27.7
$ws2.Range("G54").Value = 3.1
$ws2.Range("H54").Value = 0.74765

$ws2.Range("D55").Value = 0.22995000000000002
$ws2.Range("E55").Value = 2200
$ws2.Range("F55").Value = 27.7
$ws2.Range("G55").Value = 3.1
$ws2.Range("H55").Value = 0.74765

$ws2.Range("D56").Value = 0.22995000000000002
$ws2.Range("E56").Value = 2200
$ws2.Range("F56").Value = 27.7
$ws2.Range("G56").Value = 3.1
$ws2.Range("H56").Value = 0.74765

$ws2.Range("D57").Value = 0.22995000000000002
$ws2.Range("E57").Value = 2200
$ws2.Range("F57").Value = 27.7
$ws2.Range("G57").Value = 3.1
$ws2.Range("H57").Value = 0.74765

$ws2.Range("D58").Value = 0.23652
$ws2.Range("E58").Value = 2200
$ws2.Range("F58").Value = 27.7
$ws2.Range("G58").Value = 3.1
$ws2.Range("H58").Value = 0.74765

$ws2.Range("D59").Value = 0.23652
$ws2.Range("E59").Value = 2200
$ws2.Range("F59").Value = 27.7
$ws2.Range("G59").Value = 3.1
$ws2.Range("H59").Value = 0.74765

$ws2.Range("D60").Value = 0.23652
$ws2.Range("E60").Value = 2200
$ws2.Range("F60").Value = 27.7
$ws2.Range("G60").Value = 3.1
$ws2.Range("H60").Value = 0.74765

$ws2.Range("D61").Value = 0.2673
$ws2.Range("E61").Value = 1923
$ws2.Range("F61").Value = 22.9
$ws2.Range("G61").Value = 2.88
$ws2.Range("H61").Value = 0.8084499999999999

$ws2.Range("D62").Value = 0.2673
$ws2.Range("E62").Value = 1923
$ws2.Range("F62").Value = 22.9
$ws2.Range("G62").Value = 2.88
$ws2.Range("H62").Value = 0.8084499999999999

$ws2.Range("D63").Value = 0.2673
$ws2.Range("E63").Value = 1923
$ws2.Range("F63").Value = 22.9
$ws2.Range("G63").Value = 2.88
$ws2.Range("H63").Value = 0.8084499999999999

$ws2.Range("D64").Value = 0.22995000000000002
$ws2.Range("E64").Value = 2200
$ws2.Range("F64").Value = 27.7
$ws2.Range("G64").Value = 3.1
$ws2.Range("H64").Value = 0.74765

$ws2.Range("D65").Value = 0.22995000000000002
$ws2.Range("E65").Value = 2200
$ws2.Range("F65").Value = 27.7
$ws2.Range("G65").Value = 3.1
$ws2.Range("H65").Value = 0.74765

$ws2.Range("D66").Value = 0.22995000000000002
$ws2.Range("E66").Value = 2200
$ws2.Range("F66").Value = 27.7
$ws2.Range("G66").Value = 3.1
$ws2.Range("H66").Value = 0.74765

$ws2.Range("D67").Value = 0.22995000000000002
$ws2.Range("E67").Value = 2200
$ws2.Range("F67").Value = 27.7
$ws2.Range("G67").Value = 3.1
$ws2.Range("H67").Value = 0.74765

$ws2.Range("D68").Value = 0.22995000000000002
$ws2.Range("E68").Value = 2200
$ws2.Range("F68").Value = 27.7
$ws2.Range("G68").Value = 3.1
$ws2.Range("H68").Value = 0.74765

$ws2.Range("D69").Value = 0.23652
$ws2.Range("E69").Value = 2200
$ws2.Range("F69").Value = 27.7
$ws2.Range("G69").Value = 3.1
$ws2.Range("H69").Value = 0.74765

$ws2.Range("D70").Value = 0.21236400000000008

$ws2.Range("D71").Value = 0.21236400000000008

$ws2.Range("D72").Value = 0.21236400000000008

$ws2.Range("D73").Value = 0.22995000000000002
$ws2.Range("E73").Value = 2200
$ws2.Range("F73").Value = 27.7
$ws2.Range("G73").Value = 3.1
$ws2.Range("H73").Value = 0.74765

$ws2.Range("D74").Value = 0.22995000000000002
$ws2.Range("E74").Value = 2200
$ws2.Range("F74").Value = 27.7
$ws2.Range("G74").Value = 3.1
$ws2.Range("H74").Value = 0.74765

$ws2.Range("D75").Value = 0.22995000000000002
$ws2.Range("E75").Value = 2200
$ws2.Range("F75").Value = 27.7
$ws2.Range("G75").Value = 3.1
$ws2.Range("H75").Value = 0.74765

$ws2.Range("D76").Value = 0.22995000000000002
$ws2.Range("E76").Value = 2200
$ws2.Range("F76").Value = 27.7
$ws2.Range("G76").Value = 3.1
$ws2.Range("H76").Value = 0.74765

$ws2.Range("D77").Value = 0.22995000000000002
$ws2.Range("E77").Value = 2200
$ws2.Range("F77").Value = 27.7
$ws2.Range("G77").Value = 3.1
$ws2.Range("H77").Value = 0.74765

$ws2.Range("D78").Value = 0.22995000000000002
$ws2.Range("E78").Value = 2200
$ws2.Range("F78").Value = 27.7
$ws2.Range("G78").Value = 3.1
$ws2.Range("H78").Value = 0.74765

$ws2.Range("D79").Value = 0.22995000000000002
$ws2.Range("E79").Value = 2200
$ws2.Range("F79").Value = 27.7
$ws2.Range("G79").Value = 3.1
$ws2.Range("H79").Value = 0.74765

$ws2.Range("D80").Value = 0.22995000000000002
$ws2.Range("E80").Value = 2200
$ws2.Range("F80").Value = 27.7
$ws2.Range("G80").Value = 3.1
$ws2.Range("H80").Value = 0.74765

$ws2.Range("D81").Value = 0.22995000000000002
$ws2.Range("E81").Value = 2200
$ws2.Range("F81").Value = 27.7
$ws2.Range("G81").Value = 3.1
$ws2.Range("H81").Value = 0.74765

$ws2.Range("D82").Value = 0.08208000000000004

$ws2.Range("D83").Value = 0.08755200000000002

$ws2.Range("D84").Value = 0.3156075
$ws2.Range("E84").Value = 1726
$ws2.Range("F84").Value = 20.9
$ws2.Range("G84").Value = 2.55
$ws2.Range("H84").Value = 0.8227

$ws2.Range("D85").Value = 0.3156075
$ws2.Range("E85").Value = 1726
$ws2.Range("F85").Value = 20.9
$ws2.Range("G85").Value = 2.55
$ws2.Range("H85").Value = 0.8227

$ws2.Range("D86").Value = 0.3232125
$ws2.Range("E86").Value = 1726
$ws2.Range("F86").Value = 20.9
$ws2.Range("G86").Value = 2.55
$ws2.Range("H86").Value = 0.8227

$ws2.Range("D87").Value = 0.3232125
$ws2.Range("E87").Value = 1726
$ws2.Range("F87").Value = 20.9
$ws2.Range("G87").Value = 2.55
$ws2.Range("H87").Value = 0.8227

$ws2.Range("D88").Value = 0.3232125
$ws2.Range("E88").Value = 1726
$ws2.Range("F88").Value = 20.9
$ws2.Range("G88").Value = 2.55
$ws2.Range("H88").Value = 0.8227

$ws2.Range("D89").Value = 0.3232125
$ws2.Range("E89").Value = 1726
$ws2.Range("F89").Value = 20.9
$ws2.Range("G89").Value = 2.55
$ws2.Range("H89").Value = 0.8227

$ws2.Range("D90").Value = 0.3232125
$ws2.Range("E90").Value = 1726
$ws2.Range("F90").Value = 20.9
$ws2.Range("G90").Value = 2.55
$ws2.Range("H90").Value = 0.8227

$ws2.Range("D91").Value = 0.3156075
$ws2.Range("E91").Value = 1726
$ws2.Range("F91").Value = 20.9
$ws2.Range("G91").Value = 2.55
$ws2.Range("H91").Value = 0.8227

$ws2.Range("D92").Value = 0.3156075
$ws2.Range("E92").Value = 1726
$ws2.Range("F92").Value = 20.9
$ws2.Range("G92").Value = 2.55
$ws2.Range("H92").Value = 0.8227

$ws2.Range("D93").Value = 0.3156075
$ws2.Range("E93").Value = 1726
$ws2.Range("F93").Value = 20.9
$ws2.Range("G93").Value = 2.55
$ws2.Range("H93").Value = 0.8227

$ws2.Range("D94").Value = 0.3156075
$ws2.Range("E94").Value = 1726
$ws2.Range("F94").Value = 20.9
$ws2.Range("G94").Value = 2.55
$ws2.Range("H94").Value = 0.8227

$ws2.Range("D95").Value = 0.3156075
$ws2.Range("E95").Value = 1726
$ws2.Range("F95").Value = 20.9
$ws2.Range("G95").Value = 2.55
$ws2.Range("H95").Value = 0.8227
# --- Sheet: existing_stock ---
$ws3 = $wb.Worksheets.Item("existing_stock")

$ws3.Range("F5").Value = 0.2700000000000001
$ws3.Range("F6").Value = 0.2700000000000001
$ws3.Range("F7").Value = 0.2700000000000001
$ws3.Range("F8").Value = 0.2700000000000001
$ws3.Range("F9").Value = 0.2700000000000001
$ws3.Range("F10").Value = 0.28800000000000003
$ws3.Range("F11").Value = 0.30600000000000005
$ws3.Range("F12").Value = 0.30600000000000005
$ws3.Range("F13").Value = 0.30600000000000005
$ws3.Range("F14").Value = 0.30600000000000005
$ws3.Range("F15").Value = 0.30600000000000005
$ws3.Range("F16").Value = 0.30600000000000005
$ws3.Range("F17").Value = 0.27
$ws3.Range("F18").Value = 0.27
$ws3.Range("F19").Value = 0.27
$ws3.Range("F20").Value = 0.27
$ws3.Range("F21").Value = 0.27
$ws3.Range("F22").Value = 0.27
$ws3.Range("F23").Value = 0.27
$ws3.Range("F24").Value = 0.27
$ws3.Range("F25").Value = 0.30600000000000005
$ws3.Range("F26").Value = 0.30600000000000005
$ws3.Range("F27").Value = 0.315
$ws3.Range("F28").Value = 0.315
$ws3.Range("F29").Value = 0.315
$ws3.Range("F30").Value = 0.315
$ws3.Range("F31").Value = 0.28800000000000003
$ws3.Range("F32").Value = 0.27
$ws3.Range("F33").Value = 0.27
$ws3.Range("F34").Value = 0.27
$ws3.Range("F35").Value = 0.28800000000000003
$ws3.Range("F36").Value = 0.28800000000000003
$ws3.Range("F37").Value = 0.28800000000000003
$ws3.Range("F38").Value = 0.28800000000000003
$ws3.Range("F39").Value = 0.28800000000000003
$ws3.Range("F40").Value = 0.27
$ws3.Range("F41").Value = 0.27
$ws3.Range("F42").Value = 0.28800000000000003
$ws3.Range("F43").Value = 0.28800000000000003
$ws3.Range("F44").Value = 0.28800000000000003
$ws3.Range("F45").Value = 0.28800000000000003
$ws3.Range("F46").Value = 0.28800000000000003
$ws3.Range("F47").Value = 0.315
$ws3.Range("F48").Value = 0.315
$ws3.Range("F49").Value = 0.315
$ws3.Range("F50").Value = 0.324
$ws3.Range("F51").Value = 0.324
$ws3.Range("F52").Value = 0.324
$ws3.Range("F53").Value = 0.30600000000000005
$ws3.Range("F54").Value = 0.30600000000000005
$ws3.Range("F55").Value = 0.30600000000000005
$ws3.Range("F56").Value = 0.30600000000000005
$ws3.Range("F57").Value = 0.30600000000000005
$ws3.Range("F58").Value = 0.30600000000000005
$ws3.Range("F59").Value = 0.315
$ws3.Range("F60").Value = 0.315
$ws3.Range("F61").Value = 0.315
$ws3.Range("F62").Value = 0.315
$ws3.Range("F63").Value = 0.315
$ws3.Range("F64").Value = 0.315
$ws3.Range("F65").Value = 0.324
$ws3.Range("F66").Value = 0.324
$ws3.Range("F67").Value = 0.324
$ws3.Range("F68").Value = 0.324
$ws3.Range("F69").Value = 0.324
$ws3.Range("F70").Value = 0.324
$ws3.Range("F71").Value = 0.315
$ws3.Range("F72").Value = 0.315
$ws3.Range("F73").Value = 0.315
$ws3.Range("F74").Value = 0.315
$ws3.Range("F75").Value = 0.315
$ws3.Range("F76").Value = 0.324
$ws3.Range("F77").Value = 0.30600000000000005
$ws3.Range("F78").Value = 0.30600000000000005
$ws3.Range("F79").Value = 0.30600000000000005
$ws3.Range("F80").Value = 0.315
$ws3.Range("F81").Value = 0.315
$ws3.Range("F82").Value = 0.315
$ws3.Range("F83").Value = 0.315
$ws3.Range("F84").Value = 0.315
$ws3.Range("F85").Value = 0.315
$ws3.Range("F86").Value = 0.315
$ws3.Range("F87").Value = 0.315
$ws3.Range("F88").Value = 0.315
$ws3.Range("F89").Value = 0.24300000000000002
$ws3.Range("F90").Value = 0.2700000000000001
$ws3.Range("F91").Value = 0.2700000000000001
$ws3.Range("F92").Value = 0.2700000000000001
$ws3.Range("F93").Value = 0.2700000000000001
$ws3.Range("F94").Value = 0.3735
$ws3.Range("F95").Value = 0.3735
$ws3.Range("F96").Value = 0.3825
$ws3.Range("F97").Value = 0.3825
$ws3.Range("F98").Value = 0.3825
$ws3.Range("F99").Value = 0.3825
$ws3.Range("F100").Value = 0.3825
$ws3.Range("F101").Value = 0.3735
$ws3.Range("F102").Value = 0.3735
$ws3.Range("F103").Value = 0.3735
$ws3.Range("F104").Value = 0.3735
$ws3.Range("F105").Value = 0.3735